$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$values = @(0.9298247212236502, 0.7246052329558914, 0.9892919767590101, 0.6449069017050311, 0.8488307648806226, 0.469262412304898, 1.841566075300257, 0.01774961222854259, 0.8985662428071279, 0.4581579275178352, 0.3376508542801637, 0.6850273077074358, 1.129554360817876, 0.7141903383985572, 75.51318630499273, 120.6115918251162)

for ($r = 2; $r -le 26; $r++) {
    for ($c = 2; $c -le 17; $c++) {
        $ws.Cells.Item($r, $c).Value = $values[$c - 2]
    }
}
